$d = $word.ActiveDocument

# Update the date heading (first paragraph) - unique text, safe via Find/Replace
$d.Content.Find.Execute("2024-09-18 Wednesday", $true, $false, $false, $false,
                         $false, $true, 1, $false, "2024-09-19 Thursday", 2)

# Update the practice-problem table cells in place, preserving per-cell
# formatting. The table has 20 rows (5 "data" rows with 5 cells each,
# interleaved with blank rows); only rows 1, 5, 9, 13, 17 hold text.
$t = $d.Tables.Item(1)

$rowValues = @{
    1  = @("12÷9=", "10÷7=", "19÷6=", "70÷2=", "18÷5=")
    5  = @("40÷5=", "66÷9=", "65÷7=", "56÷3=", "95÷6=")
    9  = @("65÷9=", "70÷7=", "85÷7=", "22÷9=", "65÷9=")
    13 = @("50÷7=", "71÷9=", "86÷9=", "82÷4=", "43÷4=")
    17 = @("21÷7=", "24÷8=", "80÷6=", "69÷9=", "90÷6=")
}

foreach ($rowIndex in $rowValues.Keys) {
    $values = $rowValues[$rowIndex]
    $row = $t.Rows.Item($rowIndex)
    for ($col = 1; $col -le $values.Count; $col++) {
        $cell = $row.Cells.Item($col)
        $cell.Range.Text = $values[$col - 1]
    }
}
